$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functional KPIs")

# "Facings SOS" row (row 2), "Include Others" column (E) should now be
# "Include" instead of "Exclude" - include others for SOS.
$ws.Range("E2").Value = "Include"
